# edit.ps1 - applies the "plotter v1" slide-1 revision described by the diff:
#   1. Shape id=84 (orange roundRect) is brought to the very front of the
#      z-order (i.e. moved from near the end of the shape list to right
#      after <p:grpSpPr>, in front of shape id=82).
#   2. Shape id=82 (the big magenta roundRect) is nudged: new x/y offset.
#   3. Shape id=49 text "game play state (state = 2)" -> "File reader
#      state (state = 2)".
#   4. Shape id=86 ("Essential" textbox) is resized/repositioned and its
#      text becomes "Essential(11/23)".
#   5. Shape id=87 ("I/O status:" roundRect) is removed.
#   6. Shape id=88 ("File reader state (state = 3)" textbox) is removed.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) Shape 49: retitle "game play state" text box, keep its autosized
#    height pinned to the original value (PowerPoint's spAutoFit would
#    otherwise grow the box because the new label is longer).
# ---------------------------------------------------------------------
$shp49 = Get-ShapeById $s 49
$origHeight49 = $shp49.Height
$shp49.TextFrame.TextRange.Text = "File reader state (state = 2)"
$shp49.Height = $origHeight49

# ---------------------------------------------------------------------
# 2) Shape 82: reposition (offset only, extent unchanged).
# ---------------------------------------------------------------------
$shp82 = Get-ShapeById $s 82
$shp82.Left = 98.06520080566406
$shp82.Top = 31.013858795166016

# ---------------------------------------------------------------------
# 3) Shape 84: send to back of z-order so it becomes the first shape in
#    the tree (in front of / below shape 82, depending on direction).
# ---------------------------------------------------------------------
$shp84 = Get-ShapeById $s 84
$shp84.ZOrder(1)   # msoSendToBack

# ---------------------------------------------------------------------
# 4) Shape 86: new text + new position/size.
# ---------------------------------------------------------------------
$shp86 = Get-ShapeById $s 86
$shp86.TextFrame.TextRange.Text = "Essential(11/23)"
$shp86.Left = 328.8497009277344
$shp86.Top = 215.37875366210938
$shp86.Width = 216.39031982421875
$shp86.Height = 31.504724502563477

# ---------------------------------------------------------------------
# 5) & 6) Remove shapes 87 and 88 entirely.
# ---------------------------------------------------------------------
$shp87 = Get-ShapeById $s 87
$shp87.Delete()

$shp88 = Get-ShapeById $s 88
$shp88.Delete()
